# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") "60.530.29"
$ws.Range("E2").Value = "  +0.61%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.625.14"
$ws.Range("E3").Value = "  +1.36%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
Set-TextValue $ws.Range("D5") "581.59"
$ws.Range("E5").Value = "  +2.79%  "

# Row 6
Set-TextValue $ws.Range("D6") "144.73"
$ws.Range("E6").Value = "  +2.22%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$ws.Range("E8").Value = "  +0.48%  "

# Row 9
$ws.Range("E9").Value = "  -0.69%  "

# Row 10
$ws.Range("E10").Value = "  +1.02%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.375"
$ws.Range("E11").Value = "  +2.09%  "

# Row 12
$ws.Range("E12").Value = "  +2.93%  "

# Row 13
Set-TextValue $ws.Range("D13") "3.088.66"
$ws.Range("E13").Value = "  +1.23%  "

# Row 14
Set-TextValue $ws.Range("D14") "26.38"
$ws.Range("E14").Value = "  +13.66%  "

# Row 15
Set-TextValue $ws.Range("D15") "60.519.72"
$ws.Range("E15").Value = "  +0.55%  "

# Row 16
$ws.Range("E16").Value = "  +1.64%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.627.80"
$ws.Range("E17").Value = "  +1.17%  "

# Row 18
Set-TextValue $ws.Range("D18") "11.52"
$ws.Range("E18").Value = "  +1.89%  "

# Row 19
$ws.Range("E19").Value = "  +1.74%  "

# Row 20
Set-TextValue $ws.Range("D20") "348.19"
$ws.Range("E20").Value = "  +0.80%  "

# Row 21
Set-TextValue $ws.Range("D21") "6.89"
$ws.Range("E21").Value = "  -0.33%  "

# Row 22
$ws.Range("E22").Value = "  -0.06%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.531"
$ws.Range("E23").Value = "  -0.86%  "

# Row 24
$ws.Range("E24").Value = "  +1.43%  "

# Row 25
Set-TextValue $ws.Range("D25") "0.997"
$ws.Range("E25").Value = "  +0.16%  "

# Row 26
$ws.Range("E26").Value = "  +1.57%  "

# Row 27
Set-TextValue $ws.Range("D27") "8.13"
$ws.Range("E27").Value = "  +6.58%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.04"
$ws.Range("E28").Value = "  +14.04%  "

# Row 29
Set-TextValue $ws.Range("D29") "0.0₃0798"
$ws.Range("E29").Value = "  +2.25%  "

# Row 30
Set-TextValue $ws.Range("D30") "6.63"
$ws.Range("E30").Value = "  +5.26%  "

# Row 31
Set-TextValue $ws.Range("D31") "169.94"
$ws.Range("E31").Value = "  +5.44%  "

# Row 32
$ws.Range("E32").Value = "  +0.10%  "

# Row 33
$ws.Range("E33").Value = "  +1.06%  "

# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D34") "1.06"
$ws.Range("E34").Value = "  +10.86%  "

# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D35") "4.44"
$ws.Range("E35").Value = "  +5.28%  "

# Row 36
$ws.Range("E36").Value = "  +8.97%  "

# Row 37
$ws.Range("E37").Value = "  +3.95%  "

# Row 38
Set-TextValue $ws.Range("D38") "330.13"
$ws.Range("E38").Value = "  +12.71%  "

# Row 39
Set-TextValue $ws.Range("D39") "38.88"
$ws.Range("E39").Value = "  +3.08%  "

# Row 40
$ws.Range("E40").Value = "  +5.37%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.862"
$ws.Range("E41").Value = "  +0.99%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D42") "5.14"
$ws.Range("E42").Value = "  +4.77%  "

# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D43") "133.58"
$ws.Range("E43").Value = "  -3.29%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.100"

# Row 45
Set-TextValue $ws.Range("D45") "20.13"
$ws.Range("E45").Value = "  +3.24%  "

# Row 46
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D46") "1.00"
$ws.Range("E46").Value = "  +0.31%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "20.47"
$ws.Range("E47").Value = "  +4.32%  "

# Row 48
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D48") "0.0557"
$ws.Range("E48").Value = "  +2.62%  "

# Row 49
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D49") "0.610"
$ws.Range("E49").Value = "  +1.14%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.0244"
$ws.Range("E50").Value = "  +2.38%  "

# Row 51
Set-TextValue $ws.Range("D51") "10.75"
$ws.Range("E51").Value = "  +0.89%  "
